$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.636.94'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -4.31%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.967.39'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -6.19%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '541.21'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -5.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '152.55'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -6.90%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -3.82%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.973.94'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -6.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.112'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -4.12%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.11'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -7.96%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.366'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -5.08%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.482.60'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -6.26%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.34%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '61.681.73'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -4.35%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '23.66'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -6.53%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.970.09'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -6.01%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -6.06%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -2.84%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '380.24'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -6.99%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.95'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -6.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.65'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -6.12%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.02'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -5.30%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.470'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.64%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.084.34'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.186'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -6.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.998'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.22%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0929'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -9.46%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.24'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -7.39%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.00%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -5.45%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.44'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.62%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '158.87'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.93%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.64'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -5.46%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.00'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -5.41%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -5.87%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.28'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -5.47%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.55'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -8.65%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.92'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -4.08%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.409.78'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -10.37%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '37.01'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -4.06%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.18'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -8.15%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.665'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -4.57%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0591'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -4.75%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.30%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0245'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -5.22%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.95'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -9.24%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -3.94%  '
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.63'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -8.76%  '
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'WhiteBITCoin'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '10.47'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.21%  '
